# Add a new row of data (row 40 / item 39) describing get_gsensor_direction.sh
# to the lvp15 factory-test command list, and leave the sheet scrolled/
# selected on the newly added description cell (E40), matching the
# author's final view state as closely as the object model allows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B40").Value = "get_gsensor_direction.sh"
$ws.Range("C40").Value = "./get_gsensor_direction.sh"
$ws.Range("D40").Value = "adb shell /etc/factory-test/lvp15/get_gsensor_direction.sh"
$ws.Range("E40").Value = "get the direction of the devices through gsensor"

# Scroll the view over (topLeftCell = D1) and land the selection on the
# newly-described cell, mirroring the author's final cursor position.
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E40").Select()
